$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.292.31"
$ws.Range("D3").Value = "2.081.59"
$ws.Range("E3").Value = "  +3.18%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9987"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +1.64%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4316"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08833"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "2.082.44"
$ws.Range("E13").Value = "  +3.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.738"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.677"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001125"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06629"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.312"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "30.358.86"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("D26").Value = "2.327.98"
$ws.Range("E26").Value = "  +3.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.593"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.194"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.61%  "
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.655"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +20.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.184"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.867"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.963"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02568"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06662"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.447"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2264"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6834"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.244"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.52%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9987"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6387"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.205"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.604"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.251"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.189"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.61%  "
